$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Title / heading (appears twice: Heading1 and bold text near the end)
Replace-Text "Play Hot Spin Deluxe Free: Review & Pros and Cons | 2021" `
             "Play Hot Spin Deluxe Free: Exciting Slot Game with New Features"
Replace-Text "Play Hot Spin Deluxe Free: Review & Pros and Cons | 2021" `
             "Play Hot Spin Deluxe Free: Exciting Slot Game with New Features"

# "What we like" bullet list
Replace-Text "Innovative Hot Spin Deluxe feature" `
             "Innovative new features"
Replace-Text "Hot Spin Deluxe Multiplier can go up to x10" `
             "High-quality graphics and sound effects"
Replace-Text "Medium-high volatility and RTP rate of 96.0%" `
             "Wide range of betting options"
Replace-Text "Up to 15 free spins and 7 special features" `
             "Generous bonus features and free spins"

# "What we don't like" bullet list
Replace-Text "Limited betting range between €0.20 and €20.00 per spin" `
             "Medium-high volatility may not appeal to all players"
Replace-Text "Bonus symbol appears on reels 1, 3, and 5 only" `
             "Limited number of paylines"

# SEO description (italic text near the end)
Replace-Text "Read our Hot Spin Deluxe review and learn about the new innovative features. Play Hot Spin Deluxe online slot for free in 2021 and spin the reels with up to 15 free spins." `
             "Innovative Hot Spin Deluxe offers exciting gameplay with free spins and generous bonus features. Play Now!"
